# se calcula presupuesto mensual
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MiPresupuesto")

# --- Ingresos mensuales (fixed expenses block) ---
$ws.Range("E18").Value = 200   # Bonos
$ws.Range("E20").Value = 300   # Honorarios

# --- Egresos (gastos fijos / variables) labels ---
# Shared-string table is built in first-use order, so write the labels in
# the same order the target workbook lists them (Internet, Gas, Agua, Luz,
# Celular, Propinas, comida, cigarros, Estacionamiento).
$ws.Range("D25").Value = "Internet"
$ws.Range("D26").Value = "Gas"
$ws.Range("D27").Value = "Agua"
$ws.Range("D28").Value = "Luz"
$ws.Range("D29").Value = "Celular"
$ws.Range("D43").Value = "Propinas"
$ws.Range("D30").Value = "comida"
$ws.Range("D44").Value = "cigarros"
$ws.Range("D45").Value = "Estacionamiento"

# --- Egresos (gastos fijos) amounts ---
$ws.Range("E25").Value = 389
$ws.Range("E26").Value = 200
$ws.Range("E27").Value = 238
$ws.Range("E28").Value = 100
$ws.Range("E29").Value = 200
$ws.Range("E30").Value = 500

# --- Egresos (gastos variables) amounts ---
$ws.Range("E43").Value = 30
$ws.Range("E44").Value = 40
$ws.Range("E45").Value = 40

# --- View state: scroll position + selection ---
[void]$ws.Activate()
$excel.ActiveWindow.ScrollRow = 45
$excel.ActiveWindow.ScrollColumn = 2
[void]$ws.Range("F66").Select()

# --- Page setup: scale + orientation ---
$ws.PageSetup.Zoom = 58
$ws.PageSetup.Orientation = 1
